# Scheduled market-price refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for the affected leve rows across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 13250
$ws.Range("I47").Value = 10000
$ws.Range("J47").Value = 14333.333
$ws.Range("K47").Value = 10000
$ws.Range("L47").Value = 14333.333
$ws.Range("M47").Value = -9028
$ws.Range("N47").Value = -16277.333

$ws.Range("H53").Value = 317.91666
$ws.Range("I53").Value = 332.27274
$ws.Range("J53").Value = 305.76923
$ws.Range("K53").Value = 332.27274
$ws.Range("L53").Value = 305.76923
$ws.Range("M53").Value = 304.72726
$ws.Range("N53").Value = -1579.76923

$ws.Range("H64").Value = 146057.14
$ws.Range("I64").Value = 335466.66
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 335466.66
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -335218.66
$ws.Range("N64").Value = -4496

$ws.Range("H67").Value = 146057.14
$ws.Range("I67").Value = 335466.66
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 335466.66
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -334608.66
$ws.Range("N67").Value = -5716

$ws.Range("H129").Value = 810.1539
$ws.Range("J129").Value = 881.5
$ws.Range("L129").Value = 2644.5
$ws.Range("N129").Value = -12644.5

$ws.Range("H138").Value = 4879.915
$ws.Range("I138").Value = 1852.7778
$ws.Range("J138").Value = 6758.8276
$ws.Range("K138").Value = 5558.3334
$ws.Range("L138").Value = 20276.4828
$ws.Range("M138").Value = -418.3334000000004
$ws.Range("N138").Value = -30556.4828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28111.812
$ws.Range("I32").Value = 4539.685
$ws.Range("K32").Value = 4539.685
$ws.Range("M32").Value = -4252.685

$ws.Range("H61").Value = 2321.8147
$ws.Range("I61").Value = 1135.6364
$ws.Range("J61").Value = 3137.3125
$ws.Range("K61").Value = 1135.6364
$ws.Range("L61").Value = 3137.3125
$ws.Range("M61").Value = -923.6364000000001
$ws.Range("N61").Value = -3561.3125

$ws.Range("H105").Value = 45990
$ws.Range("J105").Value = 45990
$ws.Range("L105").Value = 45990
$ws.Range("N105").Value = -52978

$ws.Range("H124").Value = 21701.4
$ws.Range("J124").Value = 21701.4
$ws.Range("L124").Value = 21701.4
$ws.Range("N124").Value = -31521.4

$ws.Range("H132").Value = 3293.7334
$ws.Range("I132").Value = 3119.52
$ws.Range("K132").Value = 9358.559999999999
$ws.Range("M132").Value = -6828.559999999999

$ws.Range("H136").Value = 2321.8147
$ws.Range("I136").Value = 1135.6364
$ws.Range("J136").Value = 3137.3125
$ws.Range("K136").Value = 3406.9092
$ws.Range("L136").Value = 9411.9375
$ws.Range("M136").Value = -856.9092000000001
$ws.Range("N136").Value = -14511.9375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 123897.78
$ws.Range("J86").Value = 2445
$ws.Range("L86").Value = 2445
$ws.Range("N86").Value = -4691

$ws.Range("H89").Value = 123897.78
$ws.Range("J89").Value = 2445
$ws.Range("L89").Value = 12225
$ws.Range("N89").Value = -23457

$ws.Range("H105").Value = 287551.56
$ws.Range("I105").Value = 252225
$ws.Range("J105").Value = 334653.66
$ws.Range("K105").Value = 252225
$ws.Range("L105").Value = 334653.66
$ws.Range("M105").Value = -250478
$ws.Range("N105").Value = -338147.66

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15978.886
$ws.Range("I31").Value = 24230.582
$ws.Range("J31").Value = 2837.2964
$ws.Range("K31").Value = 24230.582
$ws.Range("L31").Value = 2837.2964
$ws.Range("M31").Value = -23935.582
$ws.Range("N31").Value = -3427.2964

$ws.Range("H34").Value = 15978.886
$ws.Range("I34").Value = 24230.582
$ws.Range("J34").Value = 2837.2964
$ws.Range("K34").Value = 24230.582
$ws.Range("L34").Value = 2837.2964
$ws.Range("M34").Value = -24028.582
$ws.Range("N34").Value = -3241.2964

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 667.6818
$ws.Range("J113").Value = 703.6
$ws.Range("L113").Value = 2110.8
$ws.Range("N113").Value = -6450.8

$ws.Range("H131").Value = 832.67
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 842.95874
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2528.87622
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12608.87622

$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 3000
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -13200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3149.8
$ws.Range("I122").Value = 2576.6924
$ws.Range("J122").Value = 6875
$ws.Range("K122").Value = 7730.0772
$ws.Range("L122").Value = 20625
$ws.Range("M122").Value = -5280.0772
$ws.Range("N122").Value = -25525

$ws.Range("H132").Value = 2679.5264
$ws.Range("I132").Value = 2179.9
$ws.Range("J132").Value = 3234.6667
$ws.Range("K132").Value = 6539.700000000001
$ws.Range("L132").Value = 9704.000100000001
$ws.Range("M132").Value = -4009.700000000001
$ws.Range("N132").Value = -14764.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1807.5238
$ws.Range("I55").Value = 1395
$ws.Range("J55").Value = 1904.5883
$ws.Range("K55").Value = 1395
$ws.Range("L55").Value = 1904.5883
$ws.Range("M55").Value = -1222
$ws.Range("N55").Value = -2250.5883

$ws.Range("H127").Value = 42600
$ws.Range("J127").Value = 42600
$ws.Range("L127").Value = 42600
$ws.Range("N127").Value = -52520

$ws.Range("H136").Value = 1846.3572
$ws.Range("I136").Value = 1759.9
$ws.Range("K136").Value = 5279.700000000001
$ws.Range("M136").Value = -2729.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 727
$ws.Range("I113").Value = 440.22223
$ws.Range("K113").Value = 1320.66669
$ws.Range("M113").Value = 849.33331

$ws.Range("H126").Value = 2754
$ws.Range("I126").Value = 3256.6667
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 9770.000100000001
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -7300.000100000001
$ws.Range("N126").Value = -10940

$ws.Range("H136").Value = 1426.7368
$ws.Range("J136").Value = 1951.4
$ws.Range("L136").Value = 5854.200000000001
$ws.Range("N136").Value = -10954.2
